# Apply "collate_info" processing: re-order the per-class rows (columns B:F)
# for rows 2-7 on the active sheet, cycling the row contents as produced by
# the collation step, while leaving column A (Day) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for rows 2-7, columns B (Time), C (Unit), D (Classroom),
# E (Lecturer), F (Delivery Mode) after collate_info processing.
$rows = @(
    @{ Row = 2;  B = "8:00 AM to 10:00 AM"; C = "MITS5501"; D = 10; E = "Lewis"; F = "F2F" },
    @{ Row = 3;  B = "8:00 AM to 9:00 AM";  C = "MITS5503"; D = 12; E = "Mike";  F = "F2F" },
    @{ Row = 4;  B = "8:00 AM to 9:00 AM";  C = "MITS4001"; D = 1;  E = "Jim";   F = "Online" },
    @{ Row = 5;  B = "8:00 AM to 9:00 AM";  C = "MITS5507"; D = 14; E = "Sammy"; F = "Online" },
    @{ Row = 7;  B = "8:00 AM to 9:00 AM";  C = "MITS5002"; D = 7;  E = "Mitch"; F = "Online" }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $ws.Cells.Item($rowNum, 2).Value = $r.B
    $ws.Cells.Item($rowNum, 3).Value = $r.C
    $ws.Cells.Item($rowNum, 4).Value = $r.D
    $ws.Cells.Item($rowNum, 5).Value = $r.E
    $ws.Cells.Item($rowNum, 6).Value = $r.F
}
